# Update the TCC theme/title sentence and keep Word's "_GoBack" (last-edit)
# bookmark anchored at the point where the new text was typed, exactly like
# a live editing session in Word would leave it.

$d = $word.ActiveDocument

$oldText  = "Redução do desperdício de comida por excesso de produção"
$newPart1 = "Otimização de processos de gestão de refeitórios com inteligência artificial. C"
$newPart2 = "omida por excesso de produção"

$sel = $word.Selection
$found = $sel.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if ($found) {
    $startPos = $sel.Start

    # Replace the whole matched phrase with the new wording.
    $rng = $d.Range($sel.Start, $sel.End)
    $rng.Text = $newPart1 + $newPart2

    # Word always keeps a single "_GoBack" bookmark marking the most recent
    # edit location; move it here (right after "...artificial. C", where the
    # new typing was inserted) instead of leaving it at the old position.
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }
    $bmPos = $startPos + $newPart1.Length
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
